$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N1").Value = "noOfRxAntennas"
$ws.Range("N1").Select()
